$d = $word.ActiveDocument

# Locate the last bullet item in the "characterising components" list
# ("Testing the new constant voltage state") and collapse the range to
# its end so a new list item can be appended right after it.
$rng = $d.Content
$found = $rng.Find.Execute("Testing the new constant voltage state", $true, $false, $false, $false, $false,
                            $true, 1, $false, $null, 0)

if ($found) {
    $rng.Collapse(0)

    # Insert a new paragraph after the found one; it inherits the same
    # paragraph/list formatting (Listeavsnitt style, numId 1 bullet).
    $newPara = $rng.InsertParagraphAfter()

    # The newly created paragraph is now the next paragraph after the
    # "Testing the new constant voltage state" one; fill in its text.
    $newParaRange = $rng.Next(4, 1).Paragraphs.First.Range
    $newParaRange.Text = "Find the capacity of the batteries using the given charge profile."
}
